$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Workbook calc settings: calcPr concurrentCalc="0"
#    (best-effort; engine may not persist this particular attribute)
# ------------------------------------------------------------------
try { $excel.MultiThreadedCalculation.Enabled = $false } catch {}

# ------------------------------------------------------------------
# 2. Insert a new row above the current row 28 ("Misc data sources
#    processing" / lafa_calories_weight.R row) to document the new
#    Chaudhary 2018 SI read script, pushing everything below down by
#    one row.
# ------------------------------------------------------------------
$ws.Rows.Item(28).Insert()

$ws.Range("A28").Value = "Misc data sources processing"
$ws.Range("B28").Value = "Chaudhary/read_chaudh2018_si.R"
$ws.Range("C28").Value = "raw supplemental table from Chaudhary 2018 (replaces 2015 data)"
$ws.Range("D28").Value = "raw_data/biodiversity/chaudhary2015SI/chaudhary_brooks_2018_si.xlsx"
$ws.Range("E28").Value = "processed table of ecoregion summary info and of characterization factors (CFs)"
$ws.Range("F28").Value = "raw_data/biodiversity/chaudhary2015SI/chaud2018SI_ecoregions.csv and chaud2018SI_CFs.csv"

# The row-insert carries the 7-column formatting template across the
# whole row (A:G); column G is unused in this new row, so drop it.
$ws.Range("G28").Clear()

# Row 28 needs a taller row (60) to fit the new, longer wrapped text.
$ws.Range("A28:F28").RowHeight = 60

# ------------------------------------------------------------------
# 3. scenario_analysis/scenario_prelim_biodiv.R row (now row 46, was
#    45 before the insert above) references the updated Chaudhary
#    2018 characterization-factor CSV instead of the old 2015 one.
# ------------------------------------------------------------------
$ws.Range("D46").Value = "scenarios/landflows_tnc_x_tnc_2x2x2_factorial_provisional.csv; biodiversity/chaud2018SI_CFs.csv"

# ------------------------------------------------------------------
# 4. Update the view state to match where the editor ended up:
#    scrolled/selected near the bottom of the table.
# ------------------------------------------------------------------
$ws.Range("D47").Select()
